# Applies the "training added to BayesingNetwork" class-diagram edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: wordAndCount field moves from row 7 to row 8, row 7 becomes
#     the new _lemmatizingWords field ---
$ws.Range("D7").Value = "string _lemmatizingWords"
$ws.Range("D8").Value = "Dictionary<string , int, > Array _wordAndCount"

# --- Row 9 is vacated; its contents move down into rows 10/11 along with
#     a new method (D10) ---
$ws.Range("B9").Value = ""
$ws.Range("D9").Value = ""

$ws.Range("B10").Value = "StartUp(void)"
$ws.Range("D10").Value = "pub Dictionary<s,i> WordAndcount(void)"

$ws.Range("B11").Value = "Result(string mostViable, double percent)"
$ws.Range("D11").Value = "Priv RemoveLemmatizingWords(void)"

$ws.Range("D12").Value = "Pub AddText(string text)"
$ws.Range("D13").Value = "Pub GetTotalWods(void) : count"

# --- FileReader class header becomes bold, matching the other class
#     headers (Menu / Catagory / BayesingNetwork) ---
$ws.Range("B14").Font.Bold = $true

# --- New FileReader members describing folder creation / test data ---
$ws.Range("B18").Value = "FileObj[] getTestData()"
$ws.Range("B19").Value = "FileObj[] GetSavedBayesingNetworks()"

# --- New BayesingNetwork training method ---
$ws.Range("F10").Value = "Train()"

# --- Selection / scroll position matches the new area of interest ---
$ws.Range("F10").Select()
